# Working version of schedule checker & calendar output
# Fills in instructor / meeting-space values that were determined for a
# handful of 23FQ course sections, plus the small formatting/view tweaks
# that came along with that editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("23FQ")

# --- Meeting Space (col H) assignments -----------------------------------

# CEEGR 1050 -> BANN 371 (typed in without the usual theme color, matching
# how it ended up looking after entry)
$ws.Range("H3").Value = "BANN 371"
$ws.Range("H3").Font.ThemeColor = 1

# CEEGR 3510 -> BANN 619
$ws.Range("H10").Value = "BANN 619"

# ENSC 4870 -> LEML 122
$ws.Range("H21").Value = "LEML 122"

# --- Instructor (col G) assignments --------------------------------------

# CEEGR 3510 -> Gnanapragrasam
$ws.Range("G10").Value = "Gnanapragrasam"

# CEEGR 3310 -> Riazi
$ws.Range("G8").Value = "Riazi"

# CEEGR 1000 -> Gnanapragasam
$ws.Range("G2").Value = "Gnanapragasam"

# ENSC 4870 -> Gnanapragasam
$ws.Range("G21").Value = "Gnanapragasam"

# --- Cosmetic follow-up: widen the Instructor column now that it holds
# longer names, and leave the selection where the user finished editing.

$ws.Columns.Item(7).ColumnWidth = 18.14

$ws.Range("H2").Select()
